# Code clean up, new way to create characters
#
# Regenerates the "MOOD / EMOTION / INTENSITY / EVENT / APPLIED STRATEGY"
# sample rows with a new batch of characters/events, and extends the
# "STRATEGIES RELATED" list with Attention Deployment / Cognitive Change
# entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean sheet: rows 12-13 used to carry full A:E records
# (Remorse/Fly, Remorse/BecomeRich) that go away entirely, so clearing the
# whole used range first keeps the rebuild simple and avoids leaving stale
# cells behind.
$ws.Range("A1:H22").ClearContents()

# Header row (unchanged, re-applied for a consistent rebuild)
$ws.Range("A1").Value = "MOOD     "
$ws.Range("B1").Value = "EMOTION  "
$ws.Range("C1").Value = "INTENSITY"
$ws.Range("D1").Value = "   EVENT    "
$ws.Range("E1").Value = " APPLIED STRATEGY    "
$ws.Range("F1").Value = " PERSONALITY TRAITS "
$ws.Range("G1").Value = " STRATEGIES RELATED "
$ws.Range("H1").Value = " DOMINANT PERSONALITY "

# Data rows: Mood, Emotion, Intensity, Event, Applied Strategy
$rows = @(
    @(-0.7413855195045471, "Distress", 2.387104034423828,  "Talk",         "Attention Deployment"),
    @(0,                   "Love",     1.7927955389022827, "Hello",        "None"),
    @(0.8303518295288086,  "Love",     2.6735565662384033, "Conversation", "None"),
    @(2.3478033542633057,  "Love",     4.9165802001953125, "Hug",          "None"),
    @(1.046940803527832,   "Distress", 4.1016740798950195, "Discussion",   "None"),
    @(1.7211663722991943,  "Joy",      2.2095818519592285, "Congrat",      "None"),
    @(0.7149765491485596,  "Distress", 3.1760647296905518, "Bye",          "Attention Deployment"),
    @(0,                   "Hate",     1.706213116645813,  "Fired",        "Situation Modification"),
    @(-0.544727087020874,  "Hate",     1.7539055347442627, "Crash",        "Cognitive Change"),
    @(1.7855010032653809,  "Joy",      7.482694149017334,  "Profits",      "None")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}

# Personality traits (column F), now starting at row 12 since the two
# Remorse event rows are gone
$ws.Range("F12").Value = "Low Conscientiousness"
$ws.Range("F13").Value = "High Extraversion"
$ws.Range("F14").Value = "Low Neuroticism"
$ws.Range("F15").Value = "Low Agreeableness"
$ws.Range("F16").Value = "Low Openness"

# Strategies related (column G) - two new strategies added, last one's
# intensity relabelled from Weakly to Lightly
$ws.Range("G17").Value = "[Situation Selection, Weakly]"
$ws.Range("G18").Value = "[Situation Modification, Strongly]"
$ws.Range("G19").Value = "[Attention Deployment, Strongly]"
$ws.Range("G20").Value = "[Cognitive Change, Strongly]"
$ws.Range("G21").Value = "[Response Modulation, Lightly]"

# Dominant personality (column H, unchanged)
$ws.Range("H22").Value = "Extraversion"
